# Known issues.xlsx - add two new known issues (task tree selection bug,
# scan axis calibration stop bug) and a new "Reported by" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: new known issue (task tree selection jumps to first element).
#     Fill left-to-right first (this is the order the new shared strings
#     were first introduced in the authored workbook).
$ws.Range("A5").Value = 42024
$ws.Range("B5").Value = "low"
$ws.Range("C5").Value = "DAQLab"
$ws.Range("D5").Value = "When a task controller is dropped as a child of another task controller in the Task Tree, the selection jumps back to the first element in the Task Tree. This is because the task tree is re-assembled each time an a task controller item is dragged and dropped. This gives a slightly annoying user experience when assembing task trees."
$ws.Rows.Item(5).RowHeight = 60

# --- Row 6: new known issue (scan axis calibration General Protection Fault).
$ws.Range("A6").Value = 42024
$ws.Range("B6").Value = "medium"
$ws.Range("C6").Value = "Laser Scanning, galvo axis calibration"
$ws.Range("D6").Value = "When interrupting a scan axis calibration in progress by pressing the stop button on the UITC, often a General Protection Fault error is thrown. The scan calibration should be able to stop in a more reliable way."
$ws.Rows.Item(6).RowHeight = 45

# --- Header row: insert "Reported by" between "Issue description" and
#     "Assigned to" by shifting F1 (Assigned to) -> G1 and writing the new
#     column header into E1 (old F1 "Status" -> G1, old E1 "Assigned to" -> F1).
#     Use Copy so both value AND style (bold header look) move together, and
#     work right-to-left so we never clobber a cell before reading it.
$ws.Range("F1").Copy($ws.Range("G1"))
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("E1").Value = "Reported by"

# --- "Reported by" column values (E2 introduces the new "Adrian" string).
$ws.Range("E2").Value = "Adrian"
$ws.Range("E3").Value = "Adrian"
$ws.Range("E4").Value = "Adrian"
$ws.Range("E5").Value = "Adrian"
$ws.Range("E6").Value = "Adrian"

# --- New column widths (E widened, F/G newly added). The COM layer here
#     quantizes ColumnWidth to 1/6 character steps, so these land on the
#     nearest reachable value to the authored widths.
$ws.Columns.Item(5).ColumnWidth = 29.140625
$ws.Columns.Item(6).ColumnWidth = 20.7109375
$ws.Columns.Item(7).ColumnWidth = 48.28515625

# --- Extend column G's formatting down through the rest of the table so the
#     sheet dimension/row spans widen from F to G everywhere, matching the
#     existing B/C/E/F "center, middle, wrap" style (s=7).
$ws.Range("G7:G31").HorizontalAlignment = -4108
$ws.Range("G7:G31").VerticalAlignment = -4108
$ws.Range("G7:G31").WrapText = $true

# --- Selection moves to E6.
$ws.Range("E6").Select()
